$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 350-354 (B and C values changed) ---
$ws.Cells.Item(350, 2).Value = 104691
$ws.Cells.Item(350, 3).Value = 355

$ws.Cells.Item(351, 2).Value = 104913
$ws.Cells.Item(351, 3).Value = 222

$ws.Cells.Item(352, 2).Value = 105260
$ws.Cells.Item(352, 3).Value = 347

$ws.Cells.Item(353, 2).Value = 105738
$ws.Cells.Item(353, 3).Value = 478

$ws.Cells.Item(354, 2).Value = 106350
$ws.Cells.Item(354, 3).Value = 612

# --- Append 5 new rows (355-359) for 18-22 februari 2021 ---
# Copy the formatting of row 354 down into the 5 new rows so the new
# cells pick up the same styles (text-formatted date column, etc.)
$ws.Rows("354:354").Copy()
$ws.Rows("355:359").Insert()

$dates = @("18 februari 2021", "19 februari 2021", "20 februari 2021", "21 februari 2021", "22 februari 2021")
$totals = @(106939, 107535, 107859, 108082, 108352)
$news = @(589, 596, 324, 223, 263)

for ($i = 0; $i -lt 5; $i++) {
    $r = 355 + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $totals[$i]
    $ws.Cells.Item($r, 3).Value = $news[$i]
}

# --- Update the view: selection now covers the new last two rows ---
$ws.Range("A358:A359").Select()
